$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, "ECs", "Fgf1", "Fgfr1", "ECs", 3, 1, 0.5247423333333333, 1.574227, 0.044891155074209, 0.044891155074209, 3, 1, 7.616180666666666, 22.848542, 0.06266940401417194, 0.06266940401417194, 3.996532414114888, 35.968791727034, 0.002813301934008449, 0.002813301934008449),
    @(3, "ECs", "Fgf1", "Fgfr1", "FAPs", 3, 1, 0.5247423333333333, 1.574227, 0.044891155074209, 0.044891155074209, 3, 1, 86.98680366666667, 260.960411, 0.7157670466966058, 0.7157670466966058, 45.64565832525523, 410.810924927297, 0.03213160949026593, 0.03213160949026593),
    @(4, "ECs", "Fgf1", "Fgfr1", "M2", 3, 1, 0.5247423333333333, 1.574227, 0.044891155074209, 0.044891155074209, 3, 1, 0.9360126666666666, 2.808038, 0.007701938614251506, 0.007701938614251506, 0.4911654707362222, 4.420489236626, 0.0003457489207044027, 0.0003457489207044027),
    @(5, "ECs", "Fgf1", "Fgfr1", "sCs", 3, 1, 0.5247423333333333, 1.574227, 0.044891155074209, 0.044891155074209, 3, 1, 25.99049233333334, 77.97147700000001, 0.2138616106749707, 0.2138616106749707, 13.63831159147544, 122.744804323279, 0.009600494729230219, 0.009600494729230219),
    @(6, "FAPs", "Fgf1", "Fgfr1", "ECs", 3, 1, 4.531154666666667, 13.593464, 0.3876355191593572, 0.3876355191593572, 3, 1, 7.616180666666666, 22.848542, 0.06266940401417194, 0.06266940401417194, 34.5100925699431, 310.590833129488, 0.02429288696044105, 0.02429288696044104),
    @(7, "FAPs", "Fgf1", "Fgfr1", "FAPs", 3, 1, 4.531154666666667, 13.593464, 0.3876355191593572, 0.3876355191593572, 3, 1, 86.98680366666667, 260.960411, 0.7157670466966058, 0.7157670466966058, 394.1506613726338, 3547.355952353704, 0.2774567307433987, 0.2774567307433987),
    @(8, "FAPs", "Fgf1", "Fgfr1", "M2", 3, 1, 4.531154666666667, 13.593464, 0.3876355191593572, 0.3876355191593572, 3, 1, 0.9360126666666666, 2.808038, 0.007701938614251506, 0.007701938614251506, 4.241218162625778, 38.17096346363199, 0.002985544973268883, 0.002985544973268882),
    @(9, "FAPs", "Fgf1", "Fgfr1", "sCs", 3, 1, 4.531154666666667, 13.593464, 0.3876355191593572, 0.3876355191593572, 3, 1, 25.99049233333334, 77.97147700000001, 0.2138616106749707, 0.2138616106749707, 117.7669406251476, 1059.902465626328, 0.0829003564822486, 0.08290035648224858),
    @(10, "M2", "Fgf1", "Fgfr1", "ECs", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.006844115528380393, 0.006844115528380393, 3, 1, 7.616180666666666, 22.848542, 0.06266940401417194, 0.06266940401417194, 0.6093122244215554, 5.483810019793999, 0.0004289166411677387, 0.0004289166411677387),
    @(11, "M2", "Fgf1", "Fgfr1", "FAPs", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.006844115528380393, 0.006844115528380393, 3, 1, 86.98680366666667, 260.960411, 0.7157670466966058, 0.7157670466966058, 6.959147262541889, 62.632325362877, 0.004898792358999214, 0.004898792358999214),
    @(12, "M2", "Fgf1", "Fgfr1", "M2", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.006844115528380393, 0.006844115528380393, 3, 1, 0.9360126666666666, 2.808038, 0.007701938614251506, 0.007701938614251506, 0.07488319736288888, 0.6739487762659999, 0.0000527129576684313, 0.0000527129576684313),
    @(13, "M2", "Fgf1", "Fgfr1", "sCs", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.006844115528380393, 0.006844115528380393, 3, 1, 25.99049233333334, 77.97147700000001, 0.2138616106749707, 0.2138616106749707, 2.079300031148778, 18.713700280339, 0.001463693570545009, 0.001463693570545009),
    @(14, "sCs", "Fgf1", "Fgfr1", "ECs", 3, 1, 6.553315, 19.659945, 0.5606292102380533, 0.5606292102380533, 3, 1, 7.616180666666666, 22.848542, 0.06266940401417194, 0.06266940401417194, 49.91123100557667, 449.20107905019, 0.0351342984785547, 0.0351342984785547),
    @(15, "sCs", "Fgf1", "Fgfr1", "FAPs", 3, 1, 6.553315, 19.659945, 0.5606292102380533, 0.5606292102380533, 3, 1, 86.98680366666667, 260.960411, 0.7157670466966058, 0.7157670466966058, 570.0519252708217, 5130.467327437395, 0.401279914103942, 0.401279914103942),
    @(16, "sCs", "Fgf1", "Fgfr1", "M2", 3, 1, 6.553315, 19.659945, 0.5606292102380533, 0.5606292102380533, 3, 1, 0.9360126666666666, 2.808038, 0.007701938614251506, 0.007701938614251506, 6.133985848656667, 55.20587263791, 0.004317931762609788, 0.004317931762609788),
    @(17, "sCs", "Fgf1", "Fgfr1", "sCs", 3, 1, 6.553315, 19.659945, 0.5606292102380533, 0.5606292102380533, 3, 1, 25.99049233333334, 77.97147700000001, 0.2138616106749707, 0.2138616106749707, 170.3238832654184, 1532.914949388765, 0.1198970658929468, 0.1198970658929468)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    for ($c = 0; $c -lt 16; $c++) {
        $ws.Cells.Item($r, 5 + $c).Value = $row[5 + $c]
    }
}
